# Add a "Received" column (K) to the BOM worksheet, tracking how many of
# each part have been received so far (or "-" where not yet tracked).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("K2").Value = "Received"

# Numeric "received" counts
$ws.Range("K3").Value  = 10
$ws.Range("K11").Value = 2
$ws.Range("K12").Value = 2
$ws.Range("K14").Value = 2
$ws.Range("K15").Value = 2
$ws.Range("K20").Value = 50
$ws.Range("K31").Value = 2
$ws.Range("K32").Value = 2
$ws.Range("K37").Value = 5
$ws.Range("K54").Value = 10
$ws.Range("K65").Value = 3
$ws.Range("K96").Value = 2

# Rows with no received parts yet
$ws.Range("K4").Value  = "-"
$ws.Range("K5").Value  = "-"
$ws.Range("K16").Value = "-"
$ws.Range("K21").Value = "-"
$ws.Range("K22").Value = "-"
$ws.Range("K23").Value = "-"
$ws.Range("K24").Value = "-"
$ws.Range("K25").Value = "-"
$ws.Range("K27").Value = "-"
$ws.Range("K28").Value = "-"
$ws.Range("K48").Value = "-"
$ws.Range("K49").Value = "-"

# Row with a note about a wrong part number received
$ws.Range("K45").Value = "10, but wrong part number"

# Update view state to match latest edit location
try {
    $excel.ActiveWindow.ScrollRow = 76
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("J114").Select()

Write-Host "Received column populated"
